# ProjectIndex.xlsx edit:
#  - Index sheet: PV-101/PV-102 "PAVING PLAN" rows become ABC-101/ABC-102 "ABC PLAN" rows
#  - ABC Notes sheet: NOTE 1/NOTE 2 become CONSTRUCT CURB / CONSTRUCT SIDEWALK
#  - Active tab moves from "Index" to "ABC Notes"; selections updated per sheet

$wb = $excel.ActiveWorkbook

# --- Index sheet: update Sheet / Title columns for both data rows ---
$wsIndex = $wb.Worksheets.Item("Index")
$wsIndex.Range("A2").Value = "ABC-101"
$wsIndex.Range("C2").Value = "ABC PLAN"
$wsIndex.Range("A3").Value = "ABC-102"
$wsIndex.Range("C3").Value = "ABC PLAN"
$wsIndex.Range("C2").Select()

# --- ABC Notes sheet: update the note text for rows 2 and 3 ---
$wsAbcNotes = $wb.Worksheets.Item("ABC Notes")
$wsAbcNotes.Range("B2").Value = "CONSTRUCT CURB"
$wsAbcNotes.Range("B3").Value = "CONSTRUCT SIDEWALK"

# --- Excel Notes sheet: just moves the selection ---
$wsExcelNotes = $wb.Worksheets.Item("Excel Notes")
$wsExcelNotes.Range("B3").Select()

# --- Make "ABC Notes" the active/visible tab (was "Index") ---
$wsAbcNotes.Activate()

# Best-effort: restore workbook window geometry (not guaranteed to be
# wired through this host's COM surface, included for completeness).
$win = $excel.ActiveWindow
$win.Left = -23148
$win.Top = 9684
$win.Width = 23256
$win.Height = 13896
